# Update column F ("dSF") values on Sheet1 to match the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 1
    6  = -4
    8  = -2
    9  = 7
    10 = -1
    13 = 1
    14 = -1
    22 = 1
    26 = 0
    27 = -1
    28 = 4
    31 = 4
    32 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
